$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_di6f0.png'; M=94.04347826086956; N=83.34782608695652; O=88.69565217391303; P=46; Q=10; R=10; S=10; T=10; U=10; V=10 }
  3 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_x9w7o.png'; M=92.38888888888889; N=72.94444444444444; O=82.66666666666666; P=36; Q=10; R=10; S=10; T=10; U=10; V=10 }
  4 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_bbs77.png'; M=31.64444444444445; N=21.26666666666667; O=26.45555555555556; P=45; Q=2; R=2; S=2; T=2; U=2; V=2 }
  5 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_s2zoe.png'; M=64.71428571428571; N=44.90476190476191; O=54.80952380952381; P=42; Q=5; R=5; S=5; T=5; U=5; V=5 }
  6 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_37hgm.png'; M=70.95454545454545; N=54.77272727272727; O=62.86363636363636; P=44; Q=6; R=6; S=6; T=6; U=6; V=6 }
  7 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_rru0v.png'; M=56.45238095238095; N=39.42857142857143; O=47.94047619047619; P=42; Q=4; R=4; S=4; T=4; U=4; V=4 }
  8 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_qz292.png'; M=78.26666666666667; N=59.13333333333333; O=68.7; P=45; Q=7; R=7; S=7; T=7; U=7; V=7 }
  9 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_tv8e2.png'; M=71.93023255813954; N=50.25581395348837; O=61.09302325581395; P=43; Q=6; R=6; S=6; T=6; U=6; V=6 }
  10 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_bf8nx.png'; M=86.63414634146342; N=66.63414634146342; O=76.63414634146342; P=41; Q=9; R=9; S=9; T=8; U=9; V=8 }
  11 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_kq9s9.png'; M=62.30232558139535; N=39.97674418604651; O=51.13953488372093; P=43; Q=4; R=4; S=4; T=5; U=5; V=4 }
  12 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_eiu3c.png'; M=65.1590909090909; N=46.22727272727273; O=55.69318181818181; P=44; Q=5; R=5; S=5; T=5; U=5; V=5 }
  13 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_njhlh.png'; M=59.74418604651163; N=41.51162790697674; O=50.62790697674419; P=43; Q=4; R=4; S=4; T=4; U=4; V=4 }
  14 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_bj99b.png'; M=82.79069767441861; N=65.46511627906976; O=74.12790697674419; P=43; Q=8; R=8; S=8; T=8; U=8; V=8 }
  15 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_bg264.png'; M=87.9047619047619; N=71.5; O=79.70238095238095; P=42; Q=10; R=10; S=10; T=9; U=10; V=9 }
  16 = @{ H=$null; I=$null; J='catch'; K='f'; L='stimuli/catch_02.jpg'; M=$null; N=$null; O=$null; P=$null; Q=$null; R=$null; S=$null; T=$null; U=$null; V=$null }
  17 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_vg73h.png'; M=87.7; N=72.4; O=80.05000000000001; P=50; Q=10; R=10; S=10; T=10; U=10; V=10 }
  18 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_9oofc.png'; M=82.47619047619048; N=65.5; O=73.98809523809524; P=42; Q=8; R=8; S=8; T=8; U=8; V=8 }
  19 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_sx68r.png'; M=54; N=33.2051282051282; O=43.6025641025641; P=39; Q=3; R=3; S=3; T=3; U=4; V=3 }
  20 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_qdln8.png'; M=85.51162790697674; N=67.86046511627907; O=76.68604651162791; P=43; Q=9; R=9; S=9; T=9; U=9; V=9 }
  21 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_5nlnv.png'; M=86.1219512195122; N=69.1951219512195; O=77.65853658536585; P=41; Q=9; R=9; S=9; T=9; U=9; V=9 }
  22 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_2qhro.png'; M=81.73809523809524; N=62.73809523809524; O=72.23809523809524; P=42; Q=8; R=8; S=8; T=8; U=8; V=8 }
  23 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_r10cu.png'; M=78.52380952380952; N=56.14285714285715; O=67.33333333333333; P=42; Q=7; R=7; S=7; T=7; U=7; V=6 }
  24 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_xbtev.png'; M=13.68181818181818; N=8.568181818181818; O=11.125; P=44; Q=1; R=1; S=1; T=1; U=1; V=1 }
  25 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_9684y.png'; M=77.95744680851064; N=56.70212765957447; O=67.32978723404256; P=47; Q=7; R=7; S=7; T=6; U=7; V=7 }
  26 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_jpjeg.png'; M=90.90697674418605; N=74.3953488372093; O=82.65116279069767; P=43; Q=10; R=10; S=10; T=10; U=10; V=10 }
  27 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_tujn3.png'; M=81.4090909090909; N=62.52272727272727; O=71.9659090909091; P=44; Q=8; R=8; S=8; T=7; U=8; V=7 }
  28 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_iudc4.png'; M=73.625; N=52.275; O=62.95; P=40; Q=6; R=6; S=6; T=6; U=6; V=6 }
  29 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_3ze38.png'; M=35.53191489361702; N=28.4468085106383; O=31.98936170212766; P=47; Q=2; R=2; S=2; T=3; U=2; V=3 }
  30 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_99exx.png'; M=70.02272727272727; N=51.88636363636363; O=60.95454545454545; P=44; Q=6; R=6; S=6; T=5; U=5; V=6 }
  31 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_rg4in.png'; M=49.3695652173913; N=30.21739130434782; O=39.79347826086956; P=46; Q=3; R=3; S=3; T=3; U=3; V=3 }
  32 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_tbs4n.png'; M=78.95744680851064; N=58.97872340425532; O=68.96808510638297; P=47; Q=7; R=7; S=7; T=7; U=7; V=7 }
  33 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_pey7u.png'; M=30.34883720930232; N=20.34883720930232; O=25.34883720930232; P=43; Q=1; R=2; S=2; T=2; U=2; V=2 }
  34 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_lzz3x.png'; M=18.46341463414634; N=11.92682926829268; O=15.19512195121951; P=41; Q=1; R=1; S=1; T=1; U=1; V=1 }
  35 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_vnxft.png'; M=53.22727272727273; N=34.84090909090909; O=44.03409090909091; P=44; Q=3; R=3; S=3; T=4; U=3; V=4 }
  36 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_6wne4.png'; M=25.16279069767442; N=15; O=20.08139534883721; P=43; Q=1; R=1; S=1; T=1; U=2; V=1 }
  37 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_5tr4v.png'; M=56.86046511627907; N=39.3953488372093; O=48.12790697674419; P=43; Q=4; R=4; S=4; T=4; U=4; V=4 }
  38 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_il020.png'; M=18.85416666666667; N=16.16666666666667; O=17.51041666666667; P=48; Q=1; R=1; S=1; T=1; U=1; V=1 }
  39 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_j73b6.png'; M=21.5609756097561; N=19.90243902439024; O=20.73170731707317; P=41; Q=1; R=1; S=1; T=2; U=1; V=2 }
  40 = @{ H='living_rooms'; I='target'; J='old'; K='j'; L='stimuli/img_mdh76.png'; M=37.31914893617022; N=25.12765957446809; O=31.22340425531915; P=47; Q=2; R=2; S=2; T=2; U=3; V=2 }
  41 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_4o8l0.png'; M=46.02173913043478; N=31.45652173913043; O=38.73913043478261; P=46; Q=3; R=3; S=3; T=3; U=3; V=3 }
  42 = @{ H='living_rooms'; I=$null; J='new'; K='f'; L='stimuli/img_196rk.png'; M=86.53488372093024; N=69.46511627906976; O=78; P=43; Q=9; R=9; S=9; T=9; U=9; V=9 }
}

$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($r in $rows.Keys) {
  $rowData = $rows[$r]
  foreach ($c in $cols) {
    $addr = "$c$r"
    $ws.Range($addr).Value = $rowData[$c]
  }
}
